$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 156.57143
$ws.Range("I53").Value = 93.5
$ws.Range("J53").Value = 240.66667
$ws.Range("K53").Value = 93.5
$ws.Range("L53").Value = 240.66667
$ws.Range("M53").Value = 543.5
$ws.Range("N53").Value = -1514.66667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3267.9565
$ws.Range("I70").Value = 4846.154
$ws.Range("J70").Value = 1216.3
$ws.Range("K70").Value = 14538.462
$ws.Range("L70").Value = 3648.9
$ws.Range("M70").Value = -14268.462
$ws.Range("N70").Value = -4188.9

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 3267.9565
$ws.Range("I73").Value = 4846.154
$ws.Range("J73").Value = 1216.3
$ws.Range("K73").Value = 14538.462
$ws.Range("L73").Value = 3648.9
$ws.Range("M73").Value = -13602.462
$ws.Range("N73").Value = -5520.9

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 668.26666
$ws.Range("I125").Value = 592
$ws.Range("J125").Value = 1736
$ws.Range("K125").Value = 5328
$ws.Range("L125").Value = 15624
$ws.Range("M125").Value = -2868
$ws.Range("N125").Value = -20544

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 6492.5293
$ws.Range("I132").Value = 2338.4614
$ws.Range("J132").Value = 19993.25
$ws.Range("K132").Value = 7015.3842
$ws.Range("L132").Value = 59979.75
$ws.Range("M132").Value = -4485.3842
$ws.Range("N132").Value = -65039.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1967.0209
$ws.Range("I138").Value = 1494.2307
$ws.Range("J138").Value = 2525.7727
$ws.Range("K138").Value = 4482.6921
$ws.Range("L138").Value = 7577.3181
$ws.Range("M138").Value = 657.3078999999998
$ws.Range("N138").Value = -17857.3181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1622.25
$ws.Range("I2").Value = 1970.875
$ws.Range("J2").Value = 925
$ws.Range("K2").Value = 1970.875
$ws.Range("L2").Value = 925
$ws.Range("M2").Value = -1857.875
$ws.Range("N2").Value = -1151

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1326081.6
$ws.Range("I32").Value = 1465313
$ws.Range("J32").Value = 3383.3333
$ws.Range("K32").Value = 1465313
$ws.Range("L32").Value = 3383.3333
$ws.Range("M32").Value = -1465026
$ws.Range("N32").Value = -3957.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 10931447
$ws.Range("I74").Value = 8583984
$ws.Range("J74").Value = 15199561
$ws.Range("K74").Value = 8583984
$ws.Range("L74").Value = 15199561
$ws.Range("M74").Value = -8583110
$ws.Range("N74").Value = -15201309

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 10931447
$ws.Range("I77").Value = 8583984
$ws.Range("J77").Value = 15199561
$ws.Range("K77").Value = 42919920
$ws.Range("L77").Value = 75997805
$ws.Range("M77").Value = -42915552
$ws.Range("N77").Value = -76006541

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1622.25
$ws.Range("I116").Value = 1970.875
$ws.Range("J116").Value = 925
$ws.Range("K116").Value = 1970.875
$ws.Range("L116").Value = 925
$ws.Range("M116").Value = 323.125
$ws.Range("N116").Value = -5513

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 22455.98
$ws.Range("I132").Value = 42037.92
$ws.Range("J132").Value = 2874.04
$ws.Range("K132").Value = 126113.76
$ws.Range("L132").Value = 8622.119999999999
$ws.Range("M132").Value = -123583.76
$ws.Range("N132").Value = -13682.12

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1622.25
$ws.Range("I3").Value = 1970.875
$ws.Range("J3").Value = 925
$ws.Range("K3").Value = 1970.875
$ws.Range("L3").Value = 925
$ws.Range("M3").Value = -1856.875
$ws.Range("N3").Value = -1153

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1101.1666
$ws.Range("J94").Value = 1584.6666
$ws.Range("L94").Value = 1584.6666
$ws.Range("N94").Value = -2486.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2310
$ws.Range("I105").Value = 1733.3334
$ws.Range("J105").Value = 2742.5
$ws.Range("K105").Value = 1733.3334
$ws.Range("L105").Value = 2742.5
$ws.Range("M105").Value = 13.66660000000002
$ws.Range("N105").Value = -6236.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5654.25
$ws.Range("I134").Value = 6297.3335
$ws.Range("J134").Value = 3725
$ws.Range("K134").Value = 18892.0005
$ws.Range("L134").Value = 11175
$ws.Range("M134").Value = -16357.0005
$ws.Range("N134").Value = -16245

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1906.12
$ws.Range("I31").Value = 691.81665
$ws.Range("J31").Value = 3727.575
$ws.Range("K31").Value = 691.81665
$ws.Range("L31").Value = 3727.575
$ws.Range("M31").Value = -396.81665
$ws.Range("N31").Value = -4317.575

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1906.12
$ws.Range("I34").Value = 691.81665
$ws.Range("J34").Value = 3727.575
$ws.Range("K34").Value = 691.81665
$ws.Range("L34").Value = 3727.575
$ws.Range("M34").Value = -489.81665
$ws.Range("N34").Value = -4131.575

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 9559.1
$ws.Range("I94").Value = 1665.1666
$ws.Range("K94").Value = 1665.1666
$ws.Range("M94").Value = -1214.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 64169.5
$ws.Range("I99").Value = 78439.38
$ws.Range("J99").Value = 2333.3333
$ws.Range("K99").Value = 78439.38
$ws.Range("L99").Value = 2333.3333
$ws.Range("M99").Value = -76941.38
$ws.Range("N99").Value = -5329.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 64169.5
$ws.Range("I126").Value = 78439.38
$ws.Range("J126").Value = 2333.3333
$ws.Range("K126").Value = 235318.14
$ws.Range("L126").Value = 6999.999899999999
$ws.Range("M126").Value = -232848.14
$ws.Range("N126").Value = -11939.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 9436198
$ws.Range("I132").Value = 13515175
$ws.Range("J132").Value = 3563.25
$ws.Range("K132").Value = 40545525
$ws.Range("L132").Value = 10689.75
$ws.Range("M132").Value = -40542995
$ws.Range("N132").Value = -15749.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1126.758
$ws.Range("I68").Value = 703.3889
$ws.Range("J68").Value = 1299.9546
$ws.Range("K68").Value = 2110.1667
$ws.Range("L68").Value = 3899.8638
$ws.Range("M68").Value = -1299.1667
$ws.Range("N68").Value = -5521.8638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1126.758
$ws.Range("I71").Value = 703.3889
$ws.Range("J71").Value = 1299.9546
$ws.Range("K71").Value = 6330.5001
$ws.Range("L71").Value = 11699.5914
$ws.Range("M71").Value = -2274.5001
$ws.Range("N71").Value = -19811.5914

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 942.1111
$ws.Range("I114").Value = 329.25
$ws.Range("J114").Value = 1432.4
$ws.Range("K114").Value = 987.75
$ws.Range("L114").Value = 4297.200000000001
$ws.Range("M114").Value = 2266.25
$ws.Range("N114").Value = -10805.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 3845.4546
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 3845.4546
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 11536.3638
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -18420.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1575
$ws.Range("I117").Value = 1500
$ws.Range("J117").Value = 1800
$ws.Range("K117").Value = 4500
$ws.Range("L117").Value = 5400
$ws.Range("M117").Value = -1058
$ws.Range("N117").Value = -12284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1884.8823
$ws.Range("I121").Value = 501.66666
$ws.Range("J121").Value = 2639.3635
$ws.Range("K121").Value = 1504.99998
$ws.Range("L121").Value = 7918.0905
$ws.Range("M121").Value = -194.9999800000001
$ws.Range("N121").Value = -10538.0905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1280.4615
$ws.Range("I131").Value = 1324.4445
$ws.Range("J131").Value = 1271.2559
$ws.Range("K131").Value = 3973.3335
$ws.Range("L131").Value = 3813.7677
$ws.Range("M131").Value = 1066.6665
$ws.Range("N131").Value = -13893.7677

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3986.6667
$ws.Range("I80").Value = 4265.926
$ws.Range("J80").Value = 2730
$ws.Range("K80").Value = 4265.926
$ws.Range("L80").Value = 2730
$ws.Range("M80").Value = -3267.926
$ws.Range("N80").Value = -4726

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3986.6667
$ws.Range("I83").Value = 4265.926
$ws.Range("J83").Value = 2730
$ws.Range("K83").Value = 21329.63
$ws.Range("L83").Value = 13650
$ws.Range("M83").Value = -16337.63
$ws.Range("N83").Value = -23634

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2361981.8
$ws.Range("I132").Value = 3475563.8
$ws.Range("J132").Value = 3807.8235
$ws.Range("K132").Value = 10426691.4
$ws.Range("L132").Value = 11423.4705
$ws.Range("M132").Value = -10424161.4
$ws.Range("N132").Value = -16483.4705

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1947.85
$ws.Range("I7").Value = 1893
$ws.Range("J7").Value = 2075.8333
$ws.Range("K7").Value = 1893
$ws.Range("L7").Value = 2075.8333
$ws.Range("M7").Value = -1781
$ws.Range("N7").Value = -2299.8333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1947.85
$ws.Range("I126").Value = 1893
$ws.Range("J126").Value = 2075.8333
$ws.Range("K126").Value = 5679
$ws.Range("L126").Value = 6227.499899999999
$ws.Range("M126").Value = -3209
$ws.Range("N126").Value = -11167.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7582807.5
$ws.Range("I132").Value = 2449.818
$ws.Range("J132").Value = 15163165
$ws.Range("K132").Value = 7349.454000000001
$ws.Range("L132").Value = 45489495
$ws.Range("M132").Value = -4819.454000000001
$ws.Range("N132").Value = -45494555

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3128.25
$ws.Range("I122").Value = 2806.8333
$ws.Range("J122").Value = 3449.6667
$ws.Range("K122").Value = 8420.499899999999
$ws.Range("L122").Value = 10349.0001
$ws.Range("M122").Value = -5970.499899999999
$ws.Range("N122").Value = -15249.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2070.6274
$ws.Range("I132").Value = 1386.9667
$ws.Range("J132").Value = 3047.2856
$ws.Range("K132").Value = 4160.9001
$ws.Range("L132").Value = 9141.856800000001
$ws.Range("M132").Value = -1630.9001
$ws.Range("N132").Value = -14201.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 19348250
$ws.Range("I136").Value = 21514080
$ws.Range("J136").Value = 10094255
$ws.Range("K136").Value = 64542240
$ws.Range("L136").Value = 30282765
$ws.Range("M136").Value = -64539690
$ws.Range("N136").Value = -30287865
